$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294, pushing existing rows 294-304 down to 295-305.
$ws.Rows("294:294").Insert()

# Populate the newly inserted row 294 with the new weekly record.
$ws.Range("A294").Value = 4
$ws.Range("B294").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C294").Value = "Los Lagos"
$ws.Range("D294").Value = 45041
$ws.Range("E294").Value = 10
$ws.Range("F294").Value = 100112009
$ws.Range("G294").Value = "Acelga"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 50
$ws.Range("K294").Value = 9000
$ws.Range("L294").Value = 9000
$ws.Range("M294").Value = 9000
$ws.Range("N294").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O294").Value = "Región de La Araucanía"
$ws.Range("P294").Value = 750
$ws.Range("Q294").Value = 12
$ws.Range("R294").Value = "Hortaliza"
